# "case with 380 kV done"
# Update the loading_percent results grid (A2:O25) for Case_5_76 with the
# recomputed values for the 380 kV case. Columns G, L and N stay 0 and
# column A (the index) is unchanged; only B,C,D,E,F,H,I,J,K,M,O move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 9.159246203565633
$ws.Range("C2").Value = 4.5499046280394
$ws.Range("D2").Value = 8.979792736928719
$ws.Range("E2").Value = 13.7741381894298
$ws.Range("F2").Value = 34.55325754982253
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 24.06312510769107
$ws.Range("J2").Value = 10.2682249483143
$ws.Range("K2").Value = 9.454984062778191
$ws.Range("M2").Value = 14.9799449664667
$ws.Range("O2").Value = 26.20885599705282
# Row 3
$ws.Range("B3").Value = 8.884693243138871
$ws.Range("C3").Value = 4.367182745457606
$ws.Range("D3").Value = 8.945177260699065
$ws.Range("E3").Value = 13.7707508465182
$ws.Range("F3").Value = 34.63598502883749
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 24.16142677615201
$ws.Range("J3").Value = 10.28965163120334
$ws.Range("K3").Value = 9.278017926341722
$ws.Range("M3").Value = 14.91348651162734
$ws.Range("O3").Value = 26.2980661370311
# Row 4
$ws.Range("B4").Value = 8.712933884357584
$ws.Range("C4").Value = 4.251398255476476
$ws.Range("D4").Value = 8.925247705600935
$ws.Range("E4").Value = 13.77098270969455
$ws.Range("F4").Value = 34.69421619270263
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 24.22599673055941
$ws.Range("J4").Value = 10.30397794892029
$ws.Range("K4").Value = 9.168876820606073
$ws.Range("M4").Value = 14.87470806016489
$ws.Range("O4").Value = 26.35785915870146
# Row 5
$ws.Range("B5").Value = 8.642247163269831
$ws.Range("C5").Value = 4.203379289785875
$ws.Range("D5").Value = 8.917465203559033
$ws.Range("E5").Value = 13.7716600822821
$ws.Range("F5").Value = 34.71981260068038
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 24.25336866349686
$ws.Range("J5").Value = 10.3101106098418
$ws.Range("K5").Value = 9.124336355166697
$ws.Range("M5").Value = 14.85942676693016
$ws.Range("O5").Value = 26.38348551760726
# Row 6
$ws.Range("B6").Value = 8.630471080607132
$ws.Range("C6").Value = 4.195357425595826
$ws.Range("D6").Value = 8.916193562360677
$ws.Range("E6").Value = 13.77180780581888
$ws.Range("F6").Value = 34.72417551263641
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 24.25797772076475
$ws.Range("J6").Value = 10.3111467335971
$ws.Range("K6").Value = 9.116938231947831
$ws.Range("M6").Value = 14.85692112462051
$ws.Range("O6").Value = 26.38781682421463
# Row 7
$ws.Range("B7").Value = 8.711983239181254
$ws.Range("C7").Value = 4.250753944088129
$ws.Range("D7").Value = 8.925141368112396
$ws.Range("E7").Value = 13.77098948297799
$ws.Range("F7").Value = 34.69455384088151
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 24.22636158963585
$ws.Range("J7").Value = 10.30405946302914
$ws.Range("K7").Value = 9.168276317672932
$ws.Range("M7").Value = 14.87449984612216
$ws.Range("O7").Value = 26.3581996642004
# Row 8
$ws.Range("B8").Value = 9.065302719393374
$ws.Range("C8").Value = 4.48768580520126
$ws.Range("D8").Value = 8.967586068079791
$ws.Range("E8").Value = 13.77249146527835
$ws.Range("F8").Value = 34.58023681273165
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 24.09614520729617
$ws.Range("J8").Value = 10.27537016654426
$ws.Range("K8").Value = 9.394100161046506
$ws.Range("M8").Value = 14.95661637013869
$ws.Range("O8").Value = 26.2385733496846
# Row 9
$ws.Range("B9").Value = 9.728600983315824
$ws.Range("C9").Value = 4.921230539540634
$ws.Range("D9").Value = 9.061054958811647
$ws.Range("E9").Value = 13.79370005099408
$ws.Range("F9").Value = 34.41519295800673
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 23.87421734668317
$ws.Range("J9").Value = 10.22838251832838
$ws.Range("K9").Value = 9.830679906132675
$ws.Range("M9").Value = 15.13320433210067
$ws.Range("O9").Value = 26.04386532493187
# Row 10
$ws.Range("B10").Value = 10.19269661076779
$ws.Range("C10").Value = 5.217870362417696
$ws.Range("D10").Value = 9.135577970303817
$ws.Range("E10").Value = 13.82029583442855
$ws.Range("F10").Value = 34.33014180081307
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 23.7315484329818
$ws.Range("J10").Value = 10.19949485138652
$ws.Range("K10").Value = 10.14449992273153
$ws.Range("M10").Value = 15.27168703527892
$ws.Range("O10").Value = 25.92521299511376
# Row 11
$ws.Range("B11").Value = 10.39783829902836
$ws.Range("C11").Value = 5.347575969763244
$ws.Range("D11").Value = 9.170662029413281
$ws.Range("E11").Value = 13.83475552262654
$ws.Range("F11").Value = 34.2993381287744
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 23.67107187759389
$ws.Range("J11").Value = 10.18757254545453
$ws.Range("K11").Value = 10.28514281844969
$ws.Range("M11").Value = 15.336417294183
$ws.Range("O11").Value = 25.87655130245851
# Row 12
$ws.Range("B12").Value = 10.47459216673363
$ws.Range("C12").Value = 5.395905524177963
$ws.Range("D12").Value = 9.184109639699917
$ws.Range("E12").Value = 13.84056745439154
$ws.Range("F12").Value = 34.28880910766712
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 23.64880749698178
$ws.Range("J12").Value = 10.1832328241264
$ws.Range("K12").Value = 10.33804848540862
$ws.Range("M12").Value = 15.36116346214246
$ws.Range("O12").Value = 25.85888997881394
# Row 13
$ws.Range("B13").Value = 10.45810418435092
$ws.Range("C13").Value = 5.385532388001059
$ws.Range("D13").Value = 9.181206380977567
$ws.Range("E13").Value = 13.83930084443055
$ws.Range("F13").Value = 34.29102619127629
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 23.65357419927924
$ws.Range("J13").Value = 10.18415968209164
$ws.Range("K13").Value = 10.32667073540219
$ws.Range("M13").Value = 15.35582375885731
$ws.Range("O13").Value = 25.8626595754759
# Row 14
$ws.Range("B14").Value = 10.40417187220261
$ws.Range("C14").Value = 5.35156802423229
$ws.Range("D14").Value = 9.17176517214279
$ws.Range("E14").Value = 13.83522695166331
$ws.Range("F14").Value = 34.29844913198583
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 23.66922740819249
$ws.Range("J14").Value = 10.18721200833667
$ws.Range("K14").Value = 10.28950273729186
$ws.Range("M14").Value = 15.33844856977909
$ws.Range("O14").Value = 25.87508293749696
# Row 15
$ws.Range("B15").Value = 10.37101387993498
$ws.Range("C15").Value = 5.330660435883247
$ws.Range("D15").Value = 9.166003016889308
$ws.Range("E15").Value = 13.83277528031801
$ws.Range("F15").Value = 34.30314383379194
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 23.67889839565693
$ws.Range("J15").Value = 10.18910442795898
$ws.Range("K15").Value = 10.26668893371775
$ws.Range("M15").Value = 15.3278358138411
$ws.Range("O15").Value = 25.8827923844244
# Row 16
$ws.Range("B16").Value = 10.1791637461321
$ws.Range("C16").Value = 5.209285384771841
$ws.Range("D16").Value = 9.133308266035211
$ws.Range("E16").Value = 13.81939810738206
$ws.Range("F16").Value = 34.33231373594194
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 23.7355897938511
$ws.Range("J16").Value = 10.20029850530893
$ws.Range("K16").Value = 10.13526179802729
$ws.Range("M16").Value = 15.26749036668552
$ws.Range("O16").Value = 25.92850022421335
# Row 17
$ws.Range("B17").Value = 10.05988726923605
$ws.Range("C17").Value = 5.133458499966888
$ws.Range("D17").Value = 9.113548447938866
$ws.Range("E17").Value = 13.81179422163274
$ws.Range("F17").Value = 34.35222972166928
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 23.77150168657441
$ws.Range("J17").Value = 10.20747769273148
$ws.Range("K17").Value = 10.05405828502692
$ws.Range("M17").Value = 15.23090382850442
$ws.Range("O17").Value = 25.95790272036254
# Row 18
$ws.Range("B18").Value = 9.990723605118388
$ws.Range("C18").Value = 5.089353613321949
$ws.Range("D18").Value = 9.102294989456745
$ws.Range("E18").Value = 13.80764311553937
$ws.Range("F18").Value = 34.36442711909766
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 23.79257365205588
$ws.Range("J18").Value = 10.21172171088721
$ws.Range("K18").Value = 10.00715579286722
$ws.Range("M18").Value = 15.21002433841839
$ws.Range("O18").Value = 25.97531432617533
# Row 19
$ws.Range("B19").Value = 9.967212199359867
$ws.Range("C19").Value = 5.074337199917998
$ws.Range("D19").Value = 9.098504217707267
$ws.Range("E19").Value = 13.80627591998603
$ws.Range("F19").Value = 34.36868437485315
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 23.79977975666861
$ws.Range("J19").Value = 10.21317837725175
$ws.Range("K19").Value = 9.99124321106528
$ws.Range("M19").Value = 15.20298354346155
$ws.Range("O19").Value = 25.98129543174264
# Row 20
$ws.Range("B20").Value = 10.07264283174379
$ws.Range("C20").Value = 5.141581507927325
$ws.Range("D20").Value = 9.115640394483414
$ws.Range("E20").Value = 13.81258066780867
$ws.Range("F20").Value = 34.3500327980892
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 23.76763570492163
$ws.Range("J20").Value = 10.20670158281121
$ws.Range("K20").Value = 10.06272322258928
$ws.Range("M20").Value = 15.2347816504653
$ws.Range("O20").Value = 25.9547210069602
# Row 21
$ws.Range("B21").Value = 10.42003882409992
$ws.Range("C21").Value = 5.361565790808045
$ws.Range("D21").Value = 9.17453395207604
$ws.Range("E21").Value = 13.83641445054429
$ws.Range("F21").Value = 34.29623800235515
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 23.66461239455091
$ws.Range("J21").Value = 10.18631071894342
$ws.Range("K21").Value = 10.30042982588901
$ws.Range("M21").Value = 15.34354584919393
$ws.Range("O21").Value = 25.87141309816778
# Row 22
$ws.Range("B22").Value = 10.64164109473513
$ws.Range("C22").Value = 5.500737733051984
$ws.Range("D22").Value = 9.213964863446225
$ws.Range("E22").Value = 13.85395011882936
$ws.Range("F22").Value = 34.26769992143981
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 23.60099262857394
$ws.Range("J22").Value = 10.17400399796271
$ws.Range("K22").Value = 10.45370821652944
$ws.Range("M22").Value = 15.41598754595218
$ws.Range("O22").Value = 25.82143098646829
# Row 23
$ws.Range("B23").Value = 10.52388668457107
$ws.Range("C23").Value = 5.42688998851185
$ws.Range("D23").Value = 9.192836517418773
$ws.Range("E23").Value = 13.84441287797058
$ws.Range("F23").Value = 34.28232509694357
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 23.63460785467263
$ws.Range("J23").Value = 10.18047909784655
$ws.Range("K23").Value = 10.37210598472512
$ws.Range("M23").Value = 15.37720486285875
$ws.Range("O23").Value = 25.84769834652501
# Row 24
$ws.Range("B24").Value = 10.0668778744091
$ws.Range("C24").Value = 5.137910685052312
$ws.Range("D24").Value = 9.114694292414942
$ws.Range("E24").Value = 13.81222442852593
$ws.Range("F24").Value = 34.35102369915703
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 23.76938218910164
$ws.Range("J24").Value = 10.20705209886329
$ws.Range("K24").Value = 10.05880647809282
$ws.Range("M24").Value = 15.23302800445878
$ws.Range("O24").Value = 25.95615787818876
# Row 25
$ws.Range("B25").Value = 9.552894572744206
$ws.Range("C25").Value = 4.807601966559655
$ws.Range("D25").Value = 9.034712222688961
$ws.Range("E25").Value = 13.78601866735608
$ws.Range("F25").Value = 34.45349340751814
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 23.93067593393588
$ws.Range("J25").Value = 10.22838251832838
$ws.Range("K25").Value = 9.713574018842923
$ws.Range("M25").Value = 15.08384152379505
$ws.Range("O25").Value = 26.09226038883648

Write-Output "Updated 264 cells"